$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 369
$ws1.Range("F3").Value = 2282
$ws1.Range("F6").Value = 5352
$ws1.Range("F7").Value = 382
$ws1.Range("F9").Value = 317
$ws1.Range("F10").Value = 236
$ws1.Range("F12").Value = 224
$ws1.Range("F15").Value = 131
$ws1.Range("F16").Value = 4260
$ws1.Range("F17").Value = 765
$ws1.Range("F18").Value = 780
$ws1.Range("F19").Value = 38
$ws1.Range("F23").Value = 29
$ws1.Range("G25").Value = 45
$ws1.Range("F27").Value = 17
$ws1.Range("F28").Value = 39
$ws1.Range("F29").Value = 1089
$ws1.Range("F30").Value = 12
$ws1.Range("F31").Value = 2718
$ws1.Range("F32").Value = 435
$ws1.Range("F33").Value = 187

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 369
$ws4.Range("F3").Value = 2282
$ws4.Range("F6").Value = 5352
$ws4.Range("F7").Value = 382
$ws4.Range("F9").Value = 317
$ws4.Range("F10").Value = 236
$ws4.Range("F12").Value = 224
$ws4.Range("F15").Value = 131
$ws4.Range("F16").Value = 4260
$ws4.Range("F17").Value = 765
$ws4.Range("F18").Value = 780
$ws4.Range("F19").Value = 38
$ws4.Range("F23").Value = 29
$ws4.Range("G25").Value = 45
$ws4.Range("F27").Value = 17
$ws4.Range("F29").Value = 39
$ws4.Range("F30").Value = 1089
$ws4.Range("F31").Value = 12
$ws4.Range("F32").Value = 2718
$ws4.Range("F33").Value = 435
$ws4.Range("F34").Value = 187
